$wb = $excel.ActiveWorkbook

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1572.0769
$ws.Range("J17").Value = 1574
$ws.Range("L17").Value = 4722
$ws.Range("N17").Value = -5058

$ws.Range("H33").Value = 134
$ws.Range("I33").Value = 134
$ws.Range("K33").Value = 134
$ws.Range("M33").Value = 95

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3780.4736
$ws.Range("I32").Value = 3862.8333
$ws.Range("K32").Value = 3862.8333
$ws.Range("M32").Value = -3575.8333

$ws.Range("H74").Value = 5308.45
$ws.Range("I74").Value = 3144.6
$ws.Range("K74").Value = 3144.6
$ws.Range("M74").Value = -2270.6

$ws.Range("H77").Value = 5308.45
$ws.Range("I77").Value = 3144.6
$ws.Range("K77").Value = 15723
$ws.Range("M77").Value = -11355

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5666.3335
$ws.Range("I105").Value = 7000
$ws.Range("J105").Value = 4999.5
$ws.Range("K105").Value = 7000
$ws.Range("L105").Value = 4999.5
$ws.Range("M105").Value = -5253
$ws.Range("N105").Value = -8493.5

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9022.625
$ws.Range("I31").Value = 12206
$ws.Range("K31").Value = 12206
$ws.Range("M31").Value = -11911

$ws.Range("H34").Value = 9022.625
$ws.Range("I34").Value = 12206
$ws.Range("K34").Value = 12206
$ws.Range("M34").Value = -12004

$ws.Range("H62").Value = 3100
$ws.Range("I62").Value = 1200
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 1200
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -576
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 3100
$ws.Range("I65").Value = 1200
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 6000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -2880
$ws.Range("N65").Value = -31240

$ws.Range("H86").Value = 4749.5
$ws.Range("I86").Value = 4749.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4749.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3626.5
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 4749.5
$ws.Range("I89").Value = 4749.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 23747.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -18131.5
$ws.Range("N89").ClearContents()

$ws.Range("H95").Value = 36250
$ws.Range("J95").Value = 36250
$ws.Range("L95").Value = 36250
$ws.Range("N95").Value = -41742

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 30173
$ws.Range("J34").Value = 30173
$ws.Range("L34").Value = 30173
$ws.Range("N34").Value = -30709

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H76").Value = 30173
$ws.Range("J76").Value = 30173
$ws.Range("L76").Value = 30173
$ws.Range("N76").Value = -30803

$ws.Range("H79").Value = 30173
$ws.Range("J79").Value = 30173
$ws.Range("L79").Value = 30173
$ws.Range("N79").Value = -32357

$ws.Range("H101").Value = 53500
$ws.Range("J101").Value = 53500
$ws.Range("L101").Value = 53500
$ws.Range("N101").Value = -59990

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 645.4
$ws.Range("I107").Value = 112.5
$ws.Range("J107").Value = 1000.6667
$ws.Range("K107").Value = 112.5
$ws.Range("L107").Value = 1000.6667
$ws.Range("M107").Value = 1807.5
$ws.Range("N107").Value = -4840.6667

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H113").Value = 4163.6665
$ws.Range("I113").Value = 4163.6665
$ws.Range("K113").Value = 4163.6665
$ws.Range("M113").Value = -1993.6665

$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("N114").ClearContents()

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4186.5
$ws.Range("I7").Value = 4000
$ws.Range("K7").Value = 4000
$ws.Range("M7").Value = -3888

$ws.Range("H126").Value = 4186.5
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

$ws.Range("H138").Value = 69714.5
$ws.Range("J138").Value = 69714.5
$ws.Range("L138").Value = 69714.5
$ws.Range("N138").Value = -79994.5

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 20000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H24").Value = 29999
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
